{"js": "// Office.js (Word JavaScript API) edit script.\n// Applies the same content changes described by the target diff:\n//  1. Removes \"(Supplementary Fig. 5)\" from the first paragraph.\n//  2. Removes \"(Supplementary Fig. 5h)\" from the second paragraph and adds a\n//     new parenthetical clarifying note (with the _GoBack bookmark moved here).\n//  3. Replaces \"Supplementary Figure 5f \" with \"Refer to our associated manuscript \".\n//  4. Removes the old _GoBack bookmark location (it moved above) and merges\n//     the surrounding runs.\n\nconst body = context.document.body;\n\n// --- Change 1 -----------------------------------------------------------\n// \"...the percent resistance (Supplementary Fig. 5). \" -> \"...the percent resistance. \"\nlet results = body.search(\n  \"the percent resistance (Supplementary Fig. 5). \",\n  { matchCase: true }\n);\nresults.load(\"text\");\nawait context.sync();\nif (results.items.length !== 1) {\n  throw new Error(\"Change 1: expected 1 match, found \" + results.items.length);\n}\nresults.items[0].insertText(\n  \"the percent resistance. \",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n\n// --- Change 2 -------------------------------------------------------------\n// Replace the \"(Supplementary Fig. 5h)\" parenthetical and the following\n// sentence boundary with the new clarifying text; the trailing \")\" run that\n// used to close \"(Supplementary Fig. 5h)\" is left in place (it becomes part\n// of the new \").\") text), and a new _GoBack bookmark is inserted right before\n// it.\nconst change2Find =\n  \"for each transfection (Supplementary Fig. 5h). Binding curves for the \" +\n  \"controls are used to extract the maximum observed bound donor fraction. \" +\n  \"If the difference\";\nconst change2Replace =\n  \"for each transfection. Binding curves for the controls are used to \" +\n  \"extract the maximum observed bound donor fraction (actually change in \" +\n  \"angular frequency, but referred to as bound fraction in scripts, see \" +\n  \"main text for details.). If the difference\";\n\nresults = body.search(change2Find, { matchCase: true });\nresults.load(\"text\");\nawait context.sync();\nif (results.items.length !== 1) {\n  throw new Error(\"Change 2: expected 1 match, found \" + results.items.length);\n}\nresults.items[0].insertText(change2Replace, Word.InsertLocation.replace);\nawait context.sync();\n\n// Remove the original _GoBack bookmark (it sat between \"controls. \" and\n// \"This part of the script...\" further down in the document) before\n// inserting the new one below, since a bookmark name must stay unique and\n// deleteBookmark() only removes the first occurrence found in document\n// order.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// Insert the _GoBack bookmark right after \"for details.\" and before the\n// closing \")\" that now follows it.\nresults = body.search(\"main text for details.\", { matchCase: true });\nresults.load(\"text\");\nawait context.sync();\nif (results.items.length !== 1) {\n  throw new Error(\n    \"Change 2 bookmark: expected 1 match, found \" + results.items.length\n  );\n}\nresults.items[0].getRange(Word.RangeLocation.end).insertBookmark(\"_GoBack\");\nawait context.sync();\n\n// --- Change 3 ---------------------------------------------------------\n// \"Supplementary Figure 5f shows how simulated\" ->\n// \"Refer to our associated manuscript  shows how simulated\"\n// (non-bold \"Refer to our associated manuscript\" + bold single space)\nresults = body.search(\"Supplementary Figure 5f \", { matchCase: true });\nresults.load(\"text\");\nawait context.sync();\nif (results.items.length !== 1) {\n  throw new Error(\"Change 3: expected 1 match, found \" + results.items.length);\n}\nresults.items[0].insertText(\n  \"Refer to our associated manuscript \",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n\nresults = body.search(\"Refer to our associated manuscript\", {\n  matchCase: true,\n});\nresults.load(\"text\");\nawait context.sync();\nif (results.items.length !== 1) {\n  throw new Error(\n    \"Change 3 formatting: expected 1 match, found \" + results.items.length\n  );\n}\nresults.items[0].font.bold = false;\nawait context.sync();\n\n// --- Change 4 -----------------------------------------------------------\n// The old _GoBack bookmark location (between \"controls. \" and \"This part of\n// the script...\") was already removed above (it moved to Change 2's\n// location). All that remains here is to merge the runs left behind by the\n// removed bookmark by rewriting the combined text (formatting is identical\n// on both sides, so this just collapses the run split left behind after the\n// bookmark removal).\nconst change4Text =\n  \"The third parameter used in our binding analysis, is the percent \" +\n  \"resistance. This is determined as the percent change in the observed \" +\n  \"bound fraction, at a specific free acceptor concentration, compared to \" +\n  \"the bound fraction difference between the positive and negative \" +\n  \"controls.   This part of the script is run in a semi-supervised way. \" +\n  \"The user is presented with the positive controls profile and asked if \" +\n  \"the bound fraction and the \";\nresults = body.search(change4Text, { matchCase: true });\nresults.load(\"text\");\nawait context.sync();\nif (results.items.length !== 1) {\n  throw new Error(\"Change 4: expected 1 match, found \" + results.items.length);\n}\nresults.items[0].insertText(change4Text, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# Applies the same content changes described by the target diff:\n#  1. Removes \"(Supplementary Fig. 5)\" from the first paragraph.\n#  2. Removes \"(Supplementary Fig. 5h)\" from the second paragraph and adds a\n#     new parenthetical clarifying note (with the _GoBack bookmark moved here).\n#  3. Replaces \"Supplementary Figure 5f \" with \"Refer to our associated manuscript \".\n#  4. Removes the old _GoBack bookmark location (it moved above) and merges\n#     the surrounding runs.\n\n$d = $word.ActiveDocument\n\n# --- Change 1 --------------------------------------------------------------\n# \"...the percent resistance (Supplementary Fig. 5). \" -> \"...the percent resistance. \"\n$rng = $d.Content\n$find = $rng.Find\n$find.ClearFormatting()\n$result = $find.Execute(\n    \"the percent resistance (Supplementary Fig. 5). \",\n    $true, $false, $false, $false, $false, $true, 1, $false,\n    \"the percent resistance. \", 2)\nif (-not $result) { throw \"Change 1: search text not found\" }\n\n# --- Change 2 ----------------------------------------------------------\n# Replace the \"(Supplementary Fig. 5h)\" parenthetical and the following\n# sentence boundary with the new clarifying text; the trailing \")\" that used\n# to close \"(Supplementary Fig. 5h)\" is preserved as part of the new\n# \").\" sequence that now follows the _GoBack bookmark.\n$findText2 = \"for each transfection (Supplementary Fig. 5h). Binding curves for the controls are used to extract the maximum observed bound donor fraction. If the difference\"\n$replaceText2 = \"for each transfection. Binding curves for the controls are used to extract the maximum observed bound donor fraction (actually change in angular frequency, but referred to as bound fraction in scripts, see main text for details.). If the difference\"\n\n$rng = $d.Content\n$find = $rng.Find\n$find.ClearFormatting()\n$result = $find.Execute($findText2, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText2, 2)\nif (-not $result) { throw \"Change 2: search text not found\" }\n\n# Remove the original _GoBack bookmark (it sat between \"controls. \" and\n# \"This part of the script...\" further down in the document) before adding\n# the new one below, since bookmark names must stay unique.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# Insert the _GoBack bookmark right after \"for details.\" and before the\n# closing \")\" that now follows it.\n$rng = $d.Content\n$find = $rng.Find\n$find.ClearFormatting()\n$find.Text = \"main text for details.\"\n$found = $find.Execute()\nif (-not $found) { throw \"Change 2: bookmark anchor not found\" }\n$rng.Collapse(0)  # wdCollapseEnd\n$d.Bookmarks.Add(\"_GoBack\", $rng)\n\n# --- Change 3 ----------------------------------------------------------\n# \"Supplementary Figure 5f shows how simulated\" ->\n# \"Refer to our associated manuscript  shows how simulated\"\n# (non-bold \"Refer to our associated manuscript\" + bold single space)\n$rng = $d.Content\n$find = $rng.Find\n$find.ClearFormatting()\n$result = $find.Execute(\n    \"Supplementary Figure 5f \",\n    $true, $false, $false, $false, $false, $true, 1, $false,\n    \"Refer to our associated manuscript \", 2)\nif (-not $result) { throw \"Change 3: search text not found\" }\n\n# The replacement above inherits the bold formatting of the matched text;\n# un-bold everything except the trailing space.\n$rng = $d.Content\n$find = $rng.Find\n$find.ClearFormatting()\n$find.Text = \"Refer to our associated manuscript\"\n$found = $find.Execute()\nif (-not $found) { throw \"Change 3: manuscript text not found\" }\n$rng.Font.Bold = 0\n\n# --- Change 4 ------------------------------------------------------------\n# The old _GoBack bookmark location (between \"controls. \" and \"This part of\n# the script...\") was already removed above (it moved to Change 2's\n# location). All that remains here is to merge the runs left behind by the\n# removed bookmark by rewriting the combined text (formatting is identical\n# on both sides, so this just collapses the run split left behind after the\n# bookmark removal).\n$change4Text = \"The third parameter used in our binding analysis, is the percent resistance. This is determined as the percent change in the observed bound fraction, at a specific free acceptor concentration, compared to the bound fraction difference between the positive and negative controls.   This part of the script is run in a semi-supervised way. The user is presented with the positive controls profile and asked if the bound fraction and the \"\n\n$rng = $d.Content\n$find = $rng.Find\n$find.ClearFormatting()\n$result = $find.Execute($change4Text, $true, $false, $false, $false, $false, $true, 1, $false, $change4Text, 2)\nif (-not $result) { throw \"Change 4: search text not found\" }\n"}
